$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 1: apply the "s=1" style (same as A1) to M1 and N1 (values unchanged)
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("M1:N1").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------
# B2 is removed entirely
$ws.Range("B2").Clear()

# D2, H2, I2 already carry style s=1 - just update their text
$ws.Range("D2").Value = "test"
$ws.Range("H2").Value = "test"
$ws.Range("I2").Value = "test"

# F2, K2, L2 are new cells - set value then copy format from A1 (style s=1)
$ws.Range("F2").Value = "test"
$ws.Range("K2").Value = "test"
$ws.Range("L2").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("F2").PasteSpecial($xlPasteFormats)
$ws.Range("K2").PasteSpecial($xlPasteFormats)
$ws.Range("L2").PasteSpecial($xlPasteFormats)

# M2, N2 previously had no style - set value then apply style s=1
$ws.Range("M2").Value = "test"
$ws.Range("N2").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("M2:N2").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "test"
$ws.Range("H3").Value = "test"

$ws.Range("F3").Value = "test"
$ws.Range("K3").Value = "test"
$ws.Range("L3").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("F3").PasteSpecial($xlPasteFormats)
$ws.Range("K3").PasteSpecial($xlPasteFormats)
$ws.Range("L3").PasteSpecial($xlPasteFormats)

# I3, M3, N3 previously had no style - set value then apply style s=1
$ws.Range("I3").Value = "test"
$ws.Range("M3").Value = "test"
$ws.Range("N3").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("M3:N3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "test"
$ws.Range("H4").Value = "test"

$ws.Range("F4").Value = "test"
$ws.Range("K4").Value = "test"
$ws.Range("L4").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("F4").PasteSpecial($xlPasteFormats)
$ws.Range("K4").PasteSpecial($xlPasteFormats)
$ws.Range("L4").PasteSpecial($xlPasteFormats)

# I4, M4, N4 previously had no style - set value then apply style s=1
$ws.Range("I4").Value = "test"
$ws.Range("M4").Value = "test"
$ws.Range("N4").Value = "test"
$ws.Range("A1").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("M4:N4").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

Write-Host "Edit complete"
